$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# -----------------------------------------------------------------------
# 1) Amend two existing polls: the "CAN" (national) rollup rows of poll
#    ID=1 (Leger, rows 38-43) and poll ID=2 (Leger, rows 80-85) got
#    revised sample-size (nw/nu) and vote-share (Vote) numbers.
# -----------------------------------------------------------------------

# Columns: A=ID B=Pollster C=Date D=Region E=Riding F=Party G=Vote H=nw I=nu
$amendments = @(
    # row, Vote(optional, $null = unchanged), nw, nu
    @{ Row = 38; Vote = 35;   Nw = 1238; Nu = 1242 },
    @{ Row = 39; Vote = 30;   Nw = 1238; Nu = 1242 },
    @{ Row = 40; Vote = 20;   Nw = 1238; Nu = 1242 },
    @{ Row = 41; Vote = 7;    Nw = 1238; Nu = 1242 },
    @{ Row = 42; Vote = $null; Nw = 1238; Nu = 1242 },
    @{ Row = 43; Vote = $null; Nw = 1238; Nu = 1242 },

    @{ Row = 80; Vote = 35;   Nw = 1223; Nu = 1248 },
    @{ Row = 81; Vote = 29;   Nw = 1223; Nu = 1248 },
    @{ Row = 82; Vote = 23;   Nw = 1223; Nu = 1248 },
    @{ Row = 83; Vote = 7;    Nw = 1223; Nu = 1248 },
    @{ Row = 84; Vote = 5;    Nw = 1223; Nu = 1248 },
    @{ Row = 85; Vote = $null; Nw = 1223; Nu = 1248 }
)

foreach ($a in $amendments) {
    $r = $a.Row
    if ($null -ne $a.Vote) {
        $ws.Cells.Item($r, 7).Value = $a.Vote
    }
    $ws.Cells.Item($r, 8).Value = $a.Nw
    $ws.Cells.Item($r, 9).Value = $a.Nu
}

# -----------------------------------------------------------------------
# 2) Append a brand-new poll: Abacus, date serial 44208 (2021-01-06),
#    with one row per Region/Party combination, rows 190-238.
# -----------------------------------------------------------------------

$pollId = 7
$pollster = "Abacus"
$pollDate = 44208

# Regional breakdown: each region lists its party vote shares in the
# same order the source workbook already uses (LIB, CON, NDP, BQ, GRN,
# PPC, OTH). A $null means the party cell is left blank for that region.
$regions = @(
    @{ Name = "ATL";   Votes = @{ LIB = 38; CON = 29; NDP = 17; BQ = $null; GRN = 9;  PPC = 5; OTH = 1 } },
    @{ Name = "QC";    Votes = @{ LIB = 37; CON = 14; NDP = 9;  BQ = 35;    GRN = 5;  PPC = 0; OTH = 1 } },
    @{ Name = "ON";    Votes = @{ LIB = 42; CON = 32; NDP = 17; BQ = $null; GRN = 7;  PPC = 1; OTH = 0 } },
    @{ Name = "MB/SK"; Votes = @{ LIB = 19; CON = 46; NDP = 24; BQ = $null; GRN = 3;  PPC = 5; OTH = 4 } },
    @{ Name = "AB";    Votes = @{ LIB = 21; CON = 54; NDP = 18; BQ = $null; GRN = 2;  PPC = 2; OTH = 3 } },
    @{ Name = "BC";    Votes = @{ LIB = 29; CON = 30; NDP = 29; BQ = $null; GRN = 11; PPC = 0; OTH = 0 } },
    @{ Name = "CAN";   Votes = @{ LIB = 35; CON = 31; NDP = 17; BQ = 8;     GRN = 6;  PPC = 1; OTH = 1 } }
)

# Party display order matches the order already used for every other
# poll block in the sheet.
$partyOrder = @("LIB", "CON", "NDP", "BQ", "GRN", "PPC", "OTH")

$row = 190
foreach ($region in $regions) {
    foreach ($party in $partyOrder) {
        $ws.Cells.Item($row, 1).Value = $pollId
        $ws.Cells.Item($row, 2).Value = $pollster
        $ws.Cells.Item($row, 3).Value = $pollDate
        $ws.Cells.Item($row, 3).NumberFormat = "yyyy/mm/dd"
        $ws.Cells.Item($row, 4).Value = $region.Name
        $ws.Cells.Item($row, 6).Value = $party

        $vote = $region.Votes[$party]
        if ($null -ne $vote) {
            $ws.Cells.Item($row, 7).Value = $vote
        }

        $row = $row + 1
    }
}

# -----------------------------------------------------------------------
# 3) Restore the view/selection state left behind by the author: the
#    frozen header pane scrolled down to row 47, with H85:I85 selected.
# -----------------------------------------------------------------------
[void]$ws.Range("A47").Select()
$excel.ActiveWindow.FreezePanes = $true
[void]$ws.Range("H85:I85").Select()
